$wb = $excel.ActiveWorkbook
$excel.CalculateFullRebuild()
$wb.Save()
